$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = 44165
$ws.Cells.Item(2,9).Value = 'Primera'
$ws.Cells.Item(2,10).Value = 108
$ws.Cells.Item(2,11).Value = 7000
$ws.Cells.Item(2,12).Value = 7500
$ws.Cells.Item(2,13).Value = 7222
$ws.Cells.Item(2,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(2,16).Value = 289

$ws.Cells.Item(3,4).Value = 44488
$ws.Cells.Item(3,9).Value = 'Primera'
$ws.Cells.Item(3,10).Value = 80
$ws.Cells.Item(3,11).Value = 9500
$ws.Cells.Item(3,12).Value = 10000
$ws.Cells.Item(3,13).Value = 9750
$ws.Cells.Item(3,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(3,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(3,16).Value = 390

$ws.Cells.Item(4,4).Value = 44511
$ws.Cells.Item(4,9).Value = 'Primera'
$ws.Cells.Item(4,10).Value = 80
$ws.Cells.Item(4,11).Value = 7000
$ws.Cells.Item(4,12).Value = 7500
$ws.Cells.Item(4,13).Value = 7250
$ws.Cells.Item(4,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(4,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(4,16).Value = 290

$ws.Cells.Item(5,4).Value = 44503
$ws.Cells.Item(5,9).Value = 'Primera'
$ws.Cells.Item(5,10).Value = 73
$ws.Cells.Item(5,11).Value = 7500
$ws.Cells.Item(5,12).Value = 8000
$ws.Cells.Item(5,13).Value = 7740
$ws.Cells.Item(5,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(5,16).Value = 310

$ws.Cells.Item(6,4).Value = 44410
$ws.Cells.Item(6,9).Value = 'Primera'
$ws.Cells.Item(6,10).Value = 75
$ws.Cells.Item(6,11).Value = 15000
$ws.Cells.Item(6,12).Value = 16000
$ws.Cells.Item(6,13).Value = 15533
$ws.Cells.Item(6,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(6,16).Value = 621

$ws.Cells.Item(7,4).Value = 44490
$ws.Cells.Item(7,9).Value = 'Primera'
$ws.Cells.Item(7,10).Value = 110
$ws.Cells.Item(7,11).Value = 9000
$ws.Cells.Item(7,12).Value = 9500
$ws.Cells.Item(7,13).Value = 9273
$ws.Cells.Item(7,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(7,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(7,16).Value = 371

$ws.Cells.Item(8,4).Value = 44343
$ws.Cells.Item(8,9).Value = 'Primera'
$ws.Cells.Item(8,10).Value = 18
$ws.Cells.Item(8,11).Value = 15000
$ws.Cells.Item(8,12).Value = 15000
$ws.Cells.Item(8,13).Value = 15000
$ws.Cells.Item(8,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(8,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(8,16).Value = 600

$ws.Cells.Item(9,4).Value = 44160
$ws.Cells.Item(9,9).Value = 'Primera'
$ws.Cells.Item(9,10).Value = 73
$ws.Cells.Item(9,11).Value = 7000
$ws.Cells.Item(9,12).Value = 7500
$ws.Cells.Item(9,13).Value = 7260
$ws.Cells.Item(9,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(9,16).Value = 290

$ws.Cells.Item(10,4).Value = 44365
$ws.Cells.Item(10,9).Value = 'Primera'
$ws.Cells.Item(10,10).Value = 85
$ws.Cells.Item(10,11).Value = 15000
$ws.Cells.Item(10,12).Value = 16000
$ws.Cells.Item(10,13).Value = 15529
$ws.Cells.Item(10,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(10,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(10,16).Value = 621

$ws.Cells.Item(11,4).Value = 44476
$ws.Cells.Item(11,9).Value = 'Primera'
$ws.Cells.Item(11,10).Value = 73
$ws.Cells.Item(11,11).Value = 8500
$ws.Cells.Item(11,12).Value = 9000
$ws.Cells.Item(11,13).Value = 8740
$ws.Cells.Item(11,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(11,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(11,16).Value = 350

$ws.Cells.Item(12,4).Value = 44364
$ws.Cells.Item(12,9).Value = 'Primera'
$ws.Cells.Item(12,10).Value = 85
$ws.Cells.Item(12,11).Value = 15000
$ws.Cells.Item(12,12).Value = 16000
$ws.Cells.Item(12,13).Value = 15529
$ws.Cells.Item(12,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(12,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(12,16).Value = 621

$ws.Cells.Item(13,4).Value = 44319
$ws.Cells.Item(13,9).Value = 'Primera'
$ws.Cells.Item(13,10).Value = 40
$ws.Cells.Item(13,11).Value = 15000
$ws.Cells.Item(13,12).Value = 15000
$ws.Cells.Item(13,13).Value = 15000
$ws.Cells.Item(13,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(13,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(13,16).Value = 600

$ws.Cells.Item(14,4).Value = 44414
$ws.Cells.Item(14,9).Value = 'Primera'
$ws.Cells.Item(14,10).Value = 40
$ws.Cells.Item(14,11).Value = 14000
$ws.Cells.Item(14,12).Value = 14000
$ws.Cells.Item(14,13).Value = 14000
$ws.Cells.Item(14,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(14,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(14,16).Value = 560

$ws.Cells.Item(15,4).Value = 44466
$ws.Cells.Item(15,9).Value = 'Primera'
$ws.Cells.Item(15,10).Value = 38
$ws.Cells.Item(15,11).Value = 13000
$ws.Cells.Item(15,12).Value = 13000
$ws.Cells.Item(15,13).Value = 13000
$ws.Cells.Item(15,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(15,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(15,16).Value = 520

$ws.Cells.Item(16,4).Value = 44385
$ws.Cells.Item(16,9).Value = 'Primera'
$ws.Cells.Item(16,10).Value = 35
$ws.Cells.Item(16,11).Value = 18000
$ws.Cells.Item(16,12).Value = 19000
$ws.Cells.Item(16,13).Value = 18571
$ws.Cells.Item(16,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(16,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(16,16).Value = 743

$ws.Cells.Item(17,4).Value = 44379
$ws.Cells.Item(17,9).Value = 'Primera'
$ws.Cells.Item(17,10).Value = 50
$ws.Cells.Item(17,11).Value = 16000
$ws.Cells.Item(17,12).Value = 16000
$ws.Cells.Item(17,13).Value = 16000
$ws.Cells.Item(17,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(17,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(17,16).Value = 640

$ws.Cells.Item(18,4).Value = 44484
$ws.Cells.Item(18,9).Value = 'Primera'
$ws.Cells.Item(18,10).Value = 105
$ws.Cells.Item(18,11).Value = 8000
$ws.Cells.Item(18,12).Value = 8500
$ws.Cells.Item(18,13).Value = 8238
$ws.Cells.Item(18,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(18,16).Value = 330

$ws.Cells.Item(19,4).Value = 44418
$ws.Cells.Item(19,9).Value = 'Primera'
$ws.Cells.Item(19,10).Value = 93
$ws.Cells.Item(19,11).Value = 15000
$ws.Cells.Item(19,12).Value = 16000
$ws.Cells.Item(19,13).Value = 15484
$ws.Cells.Item(19,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(19,16).Value = 619

$ws.Cells.Item(20,4).Value = 44455
$ws.Cells.Item(20,9).Value = 'Primera'
$ws.Cells.Item(20,10).Value = 70
$ws.Cells.Item(20,11).Value = 11000
$ws.Cells.Item(20,12).Value = 12000
$ws.Cells.Item(20,13).Value = 11500
$ws.Cells.Item(20,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(20,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(20,16).Value = 460

$ws.Cells.Item(21,4).Value = 44455
$ws.Cells.Item(21,9).Value = 'Segunda'
$ws.Cells.Item(21,10).Value = 38
$ws.Cells.Item(21,11).Value = 10000
$ws.Cells.Item(21,12).Value = 10000
$ws.Cells.Item(21,13).Value = 10000
$ws.Cells.Item(21,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(21,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(21,16).Value = 400

$ws.Cells.Item(22,4).Value = 44377
$ws.Cells.Item(22,9).Value = 'Primera'
$ws.Cells.Item(22,10).Value = 73
$ws.Cells.Item(22,11).Value = 16000
$ws.Cells.Item(22,12).Value = 17000
$ws.Cells.Item(22,13).Value = 16521
$ws.Cells.Item(22,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(22,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(22,16).Value = 661

$ws.Cells.Item(23,4).Value = 44397
$ws.Cells.Item(23,9).Value = 'Primera'
$ws.Cells.Item(23,10).Value = 75
$ws.Cells.Item(23,11).Value = 12000
$ws.Cells.Item(23,12).Value = 13000
$ws.Cells.Item(23,13).Value = 12533
$ws.Cells.Item(23,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(23,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(23,16).Value = 501

$ws.Cells.Item(24,4).Value = 44383
$ws.Cells.Item(24,9).Value = 'Primera'
$ws.Cells.Item(24,10).Value = 50
$ws.Cells.Item(24,11).Value = 13000
$ws.Cells.Item(24,12).Value = 13000
$ws.Cells.Item(24,13).Value = 13000
$ws.Cells.Item(24,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(24,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(24,16).Value = 520

$ws.Cells.Item(25,4).Value = 44473
$ws.Cells.Item(25,9).Value = 'Primera'
$ws.Cells.Item(25,10).Value = 85
$ws.Cells.Item(25,11).Value = 9000
$ws.Cells.Item(25,12).Value = 9500
$ws.Cells.Item(25,13).Value = 9265
$ws.Cells.Item(25,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(25,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(25,16).Value = 371

$ws.Cells.Item(26,4).Value = 44473
$ws.Cells.Item(26,9).Value = 'Segunda'
$ws.Cells.Item(26,10).Value = 40
$ws.Cells.Item(26,11).Value = 8000
$ws.Cells.Item(26,12).Value = 8000
$ws.Cells.Item(26,13).Value = 8000
$ws.Cells.Item(26,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(26,16).Value = 320

$ws.Cells.Item(27,4).Value = 44433
$ws.Cells.Item(27,9).Value = 'Primera'
$ws.Cells.Item(27,10).Value = 73
$ws.Cells.Item(27,11).Value = 14000
$ws.Cells.Item(27,12).Value = 15000
$ws.Cells.Item(27,13).Value = 14521
$ws.Cells.Item(27,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(27,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(27,16).Value = 581

$ws.Cells.Item(28,4).Value = 44159
$ws.Cells.Item(28,9).Value = 'Primera'
$ws.Cells.Item(28,10).Value = 40
$ws.Cells.Item(28,11).Value = 8000
$ws.Cells.Item(28,12).Value = 8000
$ws.Cells.Item(28,13).Value = 8000
$ws.Cells.Item(28,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(28,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(28,16).Value = 320

$ws.Cells.Item(29,4).Value = 44512
$ws.Cells.Item(29,9).Value = 'Primera'
$ws.Cells.Item(29,10).Value = 105
$ws.Cells.Item(29,11).Value = 7000
$ws.Cells.Item(29,12).Value = 7500
$ws.Cells.Item(29,13).Value = 7262
$ws.Cells.Item(29,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(29,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(29,16).Value = 290

$ws.Cells.Item(30,4).Value = 44335
$ws.Cells.Item(30,9).Value = 'Primera'
$ws.Cells.Item(30,10).Value = 35
$ws.Cells.Item(30,11).Value = 15000
$ws.Cells.Item(30,12).Value = 15000
$ws.Cells.Item(30,13).Value = 15000
$ws.Cells.Item(30,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(30,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(30,16).Value = 600

$ws.Cells.Item(31,4).Value = 44460
$ws.Cells.Item(31,9).Value = 'Primera'
$ws.Cells.Item(31,10).Value = 40
$ws.Cells.Item(31,11).Value = 11000
$ws.Cells.Item(31,12).Value = 11000
$ws.Cells.Item(31,13).Value = 11000
$ws.Cells.Item(31,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(31,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(31,16).Value = 440

$ws.Cells.Item(32,4).Value = 44489
$ws.Cells.Item(32,9).Value = 'Primera'
$ws.Cells.Item(32,10).Value = 55
$ws.Cells.Item(32,11).Value = 9000
$ws.Cells.Item(32,12).Value = 9000
$ws.Cells.Item(32,13).Value = 9000
$ws.Cells.Item(32,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(32,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(32,16).Value = 360

$ws.Cells.Item(33,4).Value = 44497
$ws.Cells.Item(33,9).Value = 'Primera'
$ws.Cells.Item(33,10).Value = 60
$ws.Cells.Item(33,11).Value = 8000
$ws.Cells.Item(33,12).Value = 8000
$ws.Cells.Item(33,13).Value = 8000
$ws.Cells.Item(33,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(33,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(33,16).Value = 320

$ws.Cells.Item(34,4).Value = 44449
$ws.Cells.Item(34,9).Value = 'Primera'
$ws.Cells.Item(34,10).Value = 55
$ws.Cells.Item(34,11).Value = 13000
$ws.Cells.Item(34,12).Value = 14000
$ws.Cells.Item(34,13).Value = 13545
$ws.Cells.Item(34,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(34,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(34,16).Value = 542

$ws.Cells.Item(35,4).Value = 44399
$ws.Cells.Item(35,9).Value = 'Primera'
$ws.Cells.Item(35,10).Value = 73
$ws.Cells.Item(35,11).Value = 1500
$ws.Cells.Item(35,12).Value = 14000
$ws.Cells.Item(35,13).Value = 7493
$ws.Cells.Item(35,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(35,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(35,16).Value = 300

$ws.Cells.Item(36,4).Value = 44482
$ws.Cells.Item(36,9).Value = 'Primera'
$ws.Cells.Item(36,10).Value = 110
$ws.Cells.Item(36,11).Value = 8000
$ws.Cells.Item(36,12).Value = 8500
$ws.Cells.Item(36,13).Value = 8273
$ws.Cells.Item(36,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(36,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(36,16).Value = 331

$ws.Cells.Item(37,4).Value = 44405
$ws.Cells.Item(37,9).Value = 'Primera'
$ws.Cells.Item(37,10).Value = 50
$ws.Cells.Item(37,11).Value = 14000
$ws.Cells.Item(37,12).Value = 14000
$ws.Cells.Item(37,13).Value = 14000
$ws.Cells.Item(37,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(37,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(37,16).Value = 560

$ws.Cells.Item(38,4).Value = 44435
$ws.Cells.Item(38,9).Value = 'Primera'
$ws.Cells.Item(38,10).Value = 73
$ws.Cells.Item(38,11).Value = 14000
$ws.Cells.Item(38,12).Value = 15000
$ws.Cells.Item(38,13).Value = 14521
$ws.Cells.Item(38,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(38,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(38,16).Value = 581

$ws.Cells.Item(39,4).Value = 44516
$ws.Cells.Item(39,9).Value = 'Primera'
$ws.Cells.Item(39,10).Value = 90
$ws.Cells.Item(39,11).Value = 7000
$ws.Cells.Item(39,12).Value = 7500
$ws.Cells.Item(39,13).Value = 7278
$ws.Cells.Item(39,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(39,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(39,16).Value = 291

$ws.Cells.Item(40,4).Value = 44168
$ws.Cells.Item(40,9).Value = 'Primera'
$ws.Cells.Item(40,10).Value = 50
$ws.Cells.Item(40,11).Value = 8500
$ws.Cells.Item(40,12).Value = 8500
$ws.Cells.Item(40,13).Value = 8500
$ws.Cells.Item(40,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(40,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(40,16).Value = 340

$ws.Cells.Item(41,4).Value = 44475
$ws.Cells.Item(41,9).Value = 'Primera'
$ws.Cells.Item(41,10).Value = 130
$ws.Cells.Item(41,11).Value = 8500
$ws.Cells.Item(41,12).Value = 9000
$ws.Cells.Item(41,13).Value = 8769
$ws.Cells.Item(41,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(41,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(41,16).Value = 351

$ws.Cells.Item(42,4).Value = 44483
$ws.Cells.Item(42,9).Value = 'Primera'
$ws.Cells.Item(42,10).Value = 90
$ws.Cells.Item(42,11).Value = 8000
$ws.Cells.Item(42,12).Value = 8500
$ws.Cells.Item(42,13).Value = 8278
$ws.Cells.Item(42,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(42,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(42,16).Value = 331

$ws.Cells.Item(43,4).Value = 44162
$ws.Cells.Item(43,9).Value = 'Primera'
$ws.Cells.Item(43,10).Value = 35
$ws.Cells.Item(43,11).Value = 7500
$ws.Cells.Item(43,12).Value = 7500
$ws.Cells.Item(43,13).Value = 7500
$ws.Cells.Item(43,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(43,15).Value = 'Provincia de Petorca'
$ws.Cells.Item(43,16).Value = 300

$ws.Cells.Item(44,4).Value = 44333
$ws.Cells.Item(44,9).Value = 'Primera'
$ws.Cells.Item(44,10).Value = 35
$ws.Cells.Item(44,11).Value = 15000
$ws.Cells.Item(44,12).Value = 15000
$ws.Cells.Item(44,13).Value = 15000
$ws.Cells.Item(44,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(44,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(44,16).Value = 600

$ws.Cells.Item(45,4).Value = 44320
$ws.Cells.Item(45,9).Value = 'Primera'
$ws.Cells.Item(45,10).Value = 60
$ws.Cells.Item(45,11).Value = 15000
$ws.Cells.Item(45,12).Value = 15000
$ws.Cells.Item(45,13).Value = 15000
$ws.Cells.Item(45,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(45,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(45,16).Value = 600

$ws.Cells.Item(46,4).Value = 44467
$ws.Cells.Item(46,9).Value = 'Primera'
$ws.Cells.Item(46,10).Value = 40
$ws.Cells.Item(46,11).Value = 14000
$ws.Cells.Item(46,12).Value = 14000
$ws.Cells.Item(46,13).Value = 14000
$ws.Cells.Item(46,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(46,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(46,16).Value = 560

$ws.Cells.Item(47,4).Value = 44390
$ws.Cells.Item(47,9).Value = 'Primera'
$ws.Cells.Item(47,10).Value = 85
$ws.Cells.Item(47,11).Value = 14000
$ws.Cells.Item(47,12).Value = 15000
$ws.Cells.Item(47,13).Value = 14471
$ws.Cells.Item(47,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(47,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(47,16).Value = 579

$ws.Cells.Item(48,4).Value = 44496
$ws.Cells.Item(48,9).Value = 'Primera'
$ws.Cells.Item(48,10).Value = 75
$ws.Cells.Item(48,11).Value = 8500
$ws.Cells.Item(48,12).Value = 9000
$ws.Cells.Item(48,13).Value = 8733
$ws.Cells.Item(48,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(48,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(48,16).Value = 349

$ws.Cells.Item(49,4).Value = 44326
$ws.Cells.Item(49,9).Value = 'Primera'
$ws.Cells.Item(49,10).Value = 35
$ws.Cells.Item(49,11).Value = 15000
$ws.Cells.Item(49,12).Value = 15000
$ws.Cells.Item(49,13).Value = 15000
$ws.Cells.Item(49,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(49,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(49,16).Value = 600

$ws.Cells.Item(50,4).Value = 44498
$ws.Cells.Item(50,9).Value = 'Primera'
$ws.Cells.Item(50,10).Value = 40
$ws.Cells.Item(50,11).Value = 8000
$ws.Cells.Item(50,12).Value = 8000
$ws.Cells.Item(50,13).Value = 8000
$ws.Cells.Item(50,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(50,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(50,16).Value = 320

$ws.Cells.Item(51,4).Value = 44321
$ws.Cells.Item(51,9).Value = 'Primera'
$ws.Cells.Item(51,10).Value = 35
$ws.Cells.Item(51,11).Value = 15000
$ws.Cells.Item(51,12).Value = 15000
$ws.Cells.Item(51,13).Value = 15000
$ws.Cells.Item(51,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(51,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(51,16).Value = 600

$ws.Cells.Item(52,4).Value = 44477
$ws.Cells.Item(52,9).Value = 'Primera'
$ws.Cells.Item(52,10).Value = 85
$ws.Cells.Item(52,11).Value = 9000
$ws.Cells.Item(52,12).Value = 10000
$ws.Cells.Item(52,13).Value = 9529
$ws.Cells.Item(52,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(52,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(52,16).Value = 381

$ws.Cells.Item(53,4).Value = 44487
$ws.Cells.Item(53,9).Value = 'Primera'
$ws.Cells.Item(53,10).Value = 100
$ws.Cells.Item(53,11).Value = 10000
$ws.Cells.Item(53,12).Value = 11000
$ws.Cells.Item(53,13).Value = 10500
$ws.Cells.Item(53,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(53,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(53,16).Value = 420

$ws.Cells.Item(54,4).Value = 44487
$ws.Cells.Item(54,9).Value = 'Segunda'
$ws.Cells.Item(54,10).Value = 45
$ws.Cells.Item(54,11).Value = 8000
$ws.Cells.Item(54,12).Value = 8000
$ws.Cells.Item(54,13).Value = 8000
$ws.Cells.Item(54,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(54,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(54,16).Value = 320

$ws.Cells.Item(55,4).Value = 44452
$ws.Cells.Item(55,9).Value = 'Primera'
$ws.Cells.Item(55,10).Value = 35
$ws.Cells.Item(55,11).Value = 15000
$ws.Cells.Item(55,12).Value = 15000
$ws.Cells.Item(55,13).Value = 15000
$ws.Cells.Item(55,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(55,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(55,16).Value = 600

$ws.Cells.Item(56,4).Value = 44505
$ws.Cells.Item(56,9).Value = 'Primera'
$ws.Cells.Item(56,10).Value = 75
$ws.Cells.Item(56,11).Value = 7000
$ws.Cells.Item(56,12).Value = 8000
$ws.Cells.Item(56,13).Value = 7467
$ws.Cells.Item(56,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(56,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(56,16).Value = 299

$ws.Cells.Item(57,4).Value = 44509
$ws.Cells.Item(57,9).Value = 'Primera'
$ws.Cells.Item(57,10).Value = 85
$ws.Cells.Item(57,11).Value = 8000
$ws.Cells.Item(57,12).Value = 8500
$ws.Cells.Item(57,13).Value = 8235
$ws.Cells.Item(57,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(57,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(57,16).Value = 329

$ws.Cells.Item(58,4).Value = 44454
$ws.Cells.Item(58,9).Value = 'Primera'
$ws.Cells.Item(58,10).Value = 35
$ws.Cells.Item(58,11).Value = 13000
$ws.Cells.Item(58,12).Value = 13000
$ws.Cells.Item(58,13).Value = 13000
$ws.Cells.Item(58,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(58,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(58,16).Value = 520

$ws.Cells.Item(59,4).Value = 44494
$ws.Cells.Item(59,9).Value = 'Primera'
$ws.Cells.Item(59,10).Value = 73
$ws.Cells.Item(59,11).Value = 8000
$ws.Cells.Item(59,12).Value = 8500
$ws.Cells.Item(59,13).Value = 8240
$ws.Cells.Item(59,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(59,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(59,16).Value = 330

$ws.Cells.Item(60,4).Value = 44356
$ws.Cells.Item(60,9).Value = 'Primera'
$ws.Cells.Item(60,10).Value = 40
$ws.Cells.Item(60,11).Value = 18000
$ws.Cells.Item(60,12).Value = 18000
$ws.Cells.Item(60,13).Value = 18000
$ws.Cells.Item(60,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(60,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(60,16).Value = 720

$ws.Cells.Item(61,4).Value = 44469
$ws.Cells.Item(61,9).Value = 'Primera'
$ws.Cells.Item(61,10).Value = 110
$ws.Cells.Item(61,11).Value = 9500
$ws.Cells.Item(61,12).Value = 10000
$ws.Cells.Item(61,13).Value = 9773
$ws.Cells.Item(61,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(61,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(61,16).Value = 391

$ws.Cells.Item(62,4).Value = 44323
$ws.Cells.Item(62,9).Value = 'Primera'
$ws.Cells.Item(62,10).Value = 58
$ws.Cells.Item(62,11).Value = 15000
$ws.Cells.Item(62,12).Value = 15000
$ws.Cells.Item(62,13).Value = 15000
$ws.Cells.Item(62,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(62,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(62,16).Value = 600

$ws.Cells.Item(63,4).Value = 44481
$ws.Cells.Item(63,9).Value = 'Primera'
$ws.Cells.Item(63,10).Value = 40
$ws.Cells.Item(63,11).Value = 9000
$ws.Cells.Item(63,12).Value = 9000
$ws.Cells.Item(63,13).Value = 9000
$ws.Cells.Item(63,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(63,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(63,16).Value = 360

$ws.Cells.Item(64,4).Value = 44417
$ws.Cells.Item(64,9).Value = 'Primera'
$ws.Cells.Item(64,10).Value = 83
$ws.Cells.Item(64,11).Value = 14000
$ws.Cells.Item(64,12).Value = 15000
$ws.Cells.Item(64,13).Value = 14542
$ws.Cells.Item(64,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(64,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(64,16).Value = 582

$ws.Cells.Item(65,4).Value = 44445
$ws.Cells.Item(65,9).Value = 'Primera'
$ws.Cells.Item(65,10).Value = 85
$ws.Cells.Item(65,11).Value = 13000
$ws.Cells.Item(65,12).Value = 14000
$ws.Cells.Item(65,13).Value = 13471
$ws.Cells.Item(65,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(65,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(65,16).Value = 539

$ws.Cells.Item(66,4).Value = 44342
$ws.Cells.Item(66,9).Value = 'Primera'
$ws.Cells.Item(66,10).Value = 35
$ws.Cells.Item(66,11).Value = 15000
$ws.Cells.Item(66,12).Value = 15000
$ws.Cells.Item(66,13).Value = 15000
$ws.Cells.Item(66,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(66,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(66,16).Value = 600

$ws.Cells.Item(67,4).Value = 44462
$ws.Cells.Item(67,9).Value = 'Primera'
$ws.Cells.Item(67,10).Value = 85
$ws.Cells.Item(67,11).Value = 11000
$ws.Cells.Item(67,12).Value = 12000
$ws.Cells.Item(67,13).Value = 11529
$ws.Cells.Item(67,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(67,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(67,16).Value = 461

$ws.Cells.Item(68,4).Value = 44406
$ws.Cells.Item(68,9).Value = 'Primera'
$ws.Cells.Item(68,10).Value = 95
$ws.Cells.Item(68,11).Value = 14000
$ws.Cells.Item(68,12).Value = 15000
$ws.Cells.Item(68,13).Value = 14474
$ws.Cells.Item(68,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(68,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(68,16).Value = 579

$ws.Cells.Item(69,4).Value = 44336
$ws.Cells.Item(69,9).Value = 'Primera'
$ws.Cells.Item(69,10).Value = 40
$ws.Cells.Item(69,11).Value = 15000
$ws.Cells.Item(69,12).Value = 15000
$ws.Cells.Item(69,13).Value = 15000
$ws.Cells.Item(69,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(69,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(69,16).Value = 600

$ws.Cells.Item(70,4).Value = 44372
$ws.Cells.Item(70,9).Value = 'Primera'
$ws.Cells.Item(70,10).Value = 55
$ws.Cells.Item(70,11).Value = 15000
$ws.Cells.Item(70,12).Value = 15000
$ws.Cells.Item(70,13).Value = 15000
$ws.Cells.Item(70,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(70,15).Value = 'Provincia de Talca'
$ws.Cells.Item(70,16).Value = 600

$ws.Cells.Item(71,4).Value = 44403
$ws.Cells.Item(71,9).Value = 'Primera'
$ws.Cells.Item(71,10).Value = 48
$ws.Cells.Item(71,11).Value = 14000
$ws.Cells.Item(71,12).Value = 14000
$ws.Cells.Item(71,13).Value = 14000
$ws.Cells.Item(71,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(71,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(71,16).Value = 560

$ws.Cells.Item(72,4).Value = 44169
$ws.Cells.Item(72,9).Value = 'Primera'
$ws.Cells.Item(72,10).Value = 38
$ws.Cells.Item(72,11).Value = 8000
$ws.Cells.Item(72,12).Value = 8000
$ws.Cells.Item(72,13).Value = 8000
$ws.Cells.Item(72,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(72,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(72,16).Value = 320

$ws.Cells.Item(73,4).Value = 44376
$ws.Cells.Item(73,9).Value = 'Primera'
$ws.Cells.Item(73,10).Value = 50
$ws.Cells.Item(73,11).Value = 16000
$ws.Cells.Item(73,12).Value = 16000
$ws.Cells.Item(73,13).Value = 16000
$ws.Cells.Item(73,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(73,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(73,16).Value = 640

$ws.Cells.Item(74,4).Value = 44474
$ws.Cells.Item(74,9).Value = 'Primera'
$ws.Cells.Item(74,10).Value = 73
$ws.Cells.Item(74,11).Value = 9000
$ws.Cells.Item(74,12).Value = 9500
$ws.Cells.Item(74,13).Value = 9260
$ws.Cells.Item(74,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(74,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(74,16).Value = 370

$ws.Cells.Item(75,4).Value = 44448
$ws.Cells.Item(75,9).Value = 'Primera'
$ws.Cells.Item(75,10).Value = 100
$ws.Cells.Item(75,11).Value = 12000
$ws.Cells.Item(75,12).Value = 13000
$ws.Cells.Item(75,13).Value = 12450
$ws.Cells.Item(75,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(75,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(75,16).Value = 498

$ws.Cells.Item(76,4).Value = 44441
$ws.Cells.Item(76,9).Value = 'Primera'
$ws.Cells.Item(76,10).Value = 82
$ws.Cells.Item(76,11).Value = 14000
$ws.Cells.Item(76,12).Value = 15000
$ws.Cells.Item(76,13).Value = 14488
$ws.Cells.Item(76,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(76,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(76,16).Value = 580

$ws.Cells.Item(77,4).Value = 44504
$ws.Cells.Item(77,9).Value = 'Primera'
$ws.Cells.Item(77,10).Value = 120
$ws.Cells.Item(77,11).Value = 7500
$ws.Cells.Item(77,12).Value = 8000
$ws.Cells.Item(77,13).Value = 7750
$ws.Cells.Item(77,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(77,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(77,16).Value = 310

$ws.Cells.Item(78,4).Value = 44350
$ws.Cells.Item(78,9).Value = 'Primera'
$ws.Cells.Item(78,10).Value = 73
$ws.Cells.Item(78,11).Value = 18000
$ws.Cells.Item(78,12).Value = 19000
$ws.Cells.Item(78,13).Value = 18521
$ws.Cells.Item(78,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(78,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(78,16).Value = 741

$ws.Cells.Item(79,4).Value = 44382
$ws.Cells.Item(79,9).Value = 'Primera'
$ws.Cells.Item(79,10).Value = 45
$ws.Cells.Item(79,11).Value = 17000
$ws.Cells.Item(79,12).Value = 17000
$ws.Cells.Item(79,13).Value = 17000
$ws.Cells.Item(79,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(79,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(79,16).Value = 680

$ws.Cells.Item(80,4).Value = 44329
$ws.Cells.Item(80,9).Value = 'Primera'
$ws.Cells.Item(80,10).Value = 35
$ws.Cells.Item(80,11).Value = 15000
$ws.Cells.Item(80,12).Value = 15000
$ws.Cells.Item(80,13).Value = 15000
$ws.Cells.Item(80,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(80,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(80,16).Value = 600

$ws.Cells.Item(81,4).Value = 44491
$ws.Cells.Item(81,9).Value = 'Primera'
$ws.Cells.Item(81,10).Value = 110
$ws.Cells.Item(81,11).Value = 8000
$ws.Cells.Item(81,12).Value = 8500
$ws.Cells.Item(81,13).Value = 8227
$ws.Cells.Item(81,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(81,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(81,16).Value = 329

$ws.Cells.Item(82,4).Value = 44166
$ws.Cells.Item(82,9).Value = 'Primera'
$ws.Cells.Item(82,10).Value = 38
$ws.Cells.Item(82,11).Value = 7000
$ws.Cells.Item(82,12).Value = 7000
$ws.Cells.Item(82,13).Value = 7000
$ws.Cells.Item(82,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(82,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(82,16).Value = 280

$ws.Cells.Item(83,4).Value = 44322
$ws.Cells.Item(83,9).Value = 'Primera'
$ws.Cells.Item(83,10).Value = 60
$ws.Cells.Item(83,11).Value = 15000
$ws.Cells.Item(83,12).Value = 15000
$ws.Cells.Item(83,13).Value = 15000
$ws.Cells.Item(83,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(83,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(83,16).Value = 600

$ws.Cells.Item(84,4).Value = 44495
$ws.Cells.Item(84,9).Value = 'Primera'
$ws.Cells.Item(84,10).Value = 80
$ws.Cells.Item(84,11).Value = 8500
$ws.Cells.Item(84,12).Value = 9000
$ws.Cells.Item(84,13).Value = 8750
$ws.Cells.Item(84,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(84,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(84,16).Value = 350

$ws.Cells.Item(85,4).Value = 44327
$ws.Cells.Item(85,9).Value = 'Primera'
$ws.Cells.Item(85,10).Value = 35
$ws.Cells.Item(85,11).Value = 15000
$ws.Cells.Item(85,12).Value = 15000
$ws.Cells.Item(85,13).Value = 15000
$ws.Cells.Item(85,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(85,15).Value = 'Provincia de Talca'
$ws.Cells.Item(85,16).Value = 600

$ws.Cells.Item(86,4).Value = 44510
$ws.Cells.Item(86,9).Value = 'Primera'
$ws.Cells.Item(86,10).Value = 90
$ws.Cells.Item(86,11).Value = 7000
$ws.Cells.Item(86,12).Value = 7500
$ws.Cells.Item(86,13).Value = 7278
$ws.Cells.Item(86,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(86,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(86,16).Value = 291

$ws.Cells.Item(87,4).Value = 44161
$ws.Cells.Item(87,9).Value = 'Primera'
$ws.Cells.Item(87,10).Value = 75
$ws.Cells.Item(87,11).Value = 7500
$ws.Cells.Item(87,12).Value = 8000
$ws.Cells.Item(87,13).Value = 7733
$ws.Cells.Item(87,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(87,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(87,16).Value = 309

$ws.Cells.Item(88,4).Value = 44468
$ws.Cells.Item(88,9).Value = 'Primera'
$ws.Cells.Item(88,10).Value = 40
$ws.Cells.Item(88,11).Value = 11000
$ws.Cells.Item(88,12).Value = 11000
$ws.Cells.Item(88,13).Value = 11000
$ws.Cells.Item(88,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(88,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(88,16).Value = 440

$ws.Cells.Item(89,4).Value = 44517
$ws.Cells.Item(89,9).Value = 'Primera'
$ws.Cells.Item(89,10).Value = 80
$ws.Cells.Item(89,11).Value = 7000
$ws.Cells.Item(89,12).Value = 7500
$ws.Cells.Item(89,13).Value = 7250
$ws.Cells.Item(89,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(89,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(89,16).Value = 290

$ws.Cells.Item(90,4).Value = 44389
$ws.Cells.Item(90,9).Value = 'Primera'
$ws.Cells.Item(90,10).Value = 97
$ws.Cells.Item(90,11).Value = 14000
$ws.Cells.Item(90,12).Value = 15000
$ws.Cells.Item(90,13).Value = 14485
$ws.Cells.Item(90,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(90,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(90,16).Value = 579

$ws.Cells.Item(91,4).Value = 44340
$ws.Cells.Item(91,9).Value = 'Primera'
$ws.Cells.Item(91,10).Value = 58
$ws.Cells.Item(91,11).Value = 16000
$ws.Cells.Item(91,12).Value = 16000
$ws.Cells.Item(91,13).Value = 16000
$ws.Cells.Item(91,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(91,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(91,16).Value = 640

$ws.Cells.Item(92,4).Value = 44340
$ws.Cells.Item(92,9).Value = 'Primera'
$ws.Cells.Item(92,10).Value = 47
$ws.Cells.Item(92,11).Value = 15000
$ws.Cells.Item(92,12).Value = 15000
$ws.Cells.Item(92,13).Value = 15000
$ws.Cells.Item(92,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(92,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(92,16).Value = 600

$ws.Cells.Item(93,4).Value = 44515
$ws.Cells.Item(93,9).Value = 'Primera'
$ws.Cells.Item(93,10).Value = 115
$ws.Cells.Item(93,11).Value = 7500
$ws.Cells.Item(93,12).Value = 8000
$ws.Cells.Item(93,13).Value = 7761
$ws.Cells.Item(93,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(93,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(93,16).Value = 310

$ws.Cells.Item(94,4).Value = 44330
$ws.Cells.Item(94,9).Value = 'Primera'
$ws.Cells.Item(94,10).Value = 40
$ws.Cells.Item(94,11).Value = 15000
$ws.Cells.Item(94,12).Value = 15000
$ws.Cells.Item(94,13).Value = 15000
$ws.Cells.Item(94,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(94,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(94,16).Value = 600

$ws.Cells.Item(95,4).Value = 44400
$ws.Cells.Item(95,9).Value = 'Primera'
$ws.Cells.Item(95,10).Value = 40
$ws.Cells.Item(95,11).Value = 15000
$ws.Cells.Item(95,12).Value = 15000
$ws.Cells.Item(95,13).Value = 15000
$ws.Cells.Item(95,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(95,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(95,16).Value = 600

$ws.Cells.Item(96,1).Value = 3
$ws.Cells.Item(96,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(96,3).Value = 'Coquimbo'
$ws.Cells.Item(96,4).Value = 44508
$ws.Cells.Item(96,5).Value = 5
$ws.Cells.Item(96,6).Value = 100112026
$ws.Cells.Item(96,7).Value = 'Haba'
$ws.Cells.Item(96,8).Value = 'Sin especificar'
$ws.Cells.Item(96,9).Value = 'Primera'
$ws.Cells.Item(96,10).Value = 90
$ws.Cells.Item(96,11).Value = 7000
$ws.Cells.Item(96,12).Value = 7500
$ws.Cells.Item(96,13).Value = 7278
$ws.Cells.Item(96,14).Value = '$/malla 25 kilos'
$ws.Cells.Item(96,15).Value = 'Provincia de Quillota'
$ws.Cells.Item(96,16).Value = 291
$ws.Cells.Item(96,17).Value = 25
$ws.Cells.Item(96,18).Value = 'Hortaliza'

$ws.Range("D96").NumberFormat = $ws.Range("D95").NumberFormat